$d = $word.ActiveDocument

# Title
$d.Content.Find.Execute("La Palma Earthquakes", $true, $false, $false, $false, $false, $true, 1, $false, "Manuscript example for a paper", 2)

# Author 1
$d.Content.Find.Execute("Steve Purves", $true, $false, $false, $false, $false, $true, 1, $false, "Chiara Fichera", 2)

# Author 2
$d.Content.Find.Execute("Rowan Cockett", $true, $false, $false, $false, $false, $true, 1, $false, "John Doe", 2)

# Abstract body
$d.Content.Find.Execute("In September 2021, a significant jump in seismic activity on the island of La Palma (Canary Islands, Spain) signaled the start of a volcanic crisis that still continues at the time of writing. Earthquake data is continually collected and published by the Instituto Geográphico Nacional (IGN). …", $true, $false, $false, $false, $false, $true, 1, $false, "This is an attempt to create a project to store text and code for a paper. Maybe this is cool, we don´t know yet. Bla bla bla…", 2)
